# Insert two new rows at row 44, pushing the existing rows 44-63 down to 46-65.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("44:45").Insert()

# Row 44: new weekly "Flame Seedless" price entry (Provincia de Limarí, $/caja 15 kilos)
$ws.Range("A44").Value = 8
$ws.Range('B44').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C44').Value = 'Coquimbo'
$ws.Range("D44").Value = 44572
$ws.Range("E44").Value = 4
$ws.Range('F44').Value = 'Fruta'
$ws.Range("G44").Value = 100109
$ws.Range('H44').Value = 'Uva'
$ws.Range("I44").Value = 100109001
$ws.Range('J44').Value = 'Uva'
$ws.Range('K44').Value = 'Flame Seedless'
$ws.Range('L44').Value = 'Primera'
$ws.Range("M44").Value = 520
$ws.Range("N44").Value = 9500
$ws.Range("O44").Value = 10000
$ws.Range("P44").Value = 9750
$ws.Range('Q44').Value = '$/caja 15 kilos'
$ws.Range('R44').Value = 'Provincia de Limarí'
$ws.Range("S44").Value = 650
$ws.Range("T44").Value = 15

# Row 45: new weekly "Superior Seedless" price entry (Provincia de Limarí, $/caja 15 kilos)
$ws.Range("A45").Value = 8
$ws.Range('B45').Value = 'Terminal La Palmera de La Serena'
$ws.Range('C45').Value = 'Coquimbo'
$ws.Range("D45").Value = 44572
$ws.Range("E45").Value = 4
$ws.Range('F45').Value = 'Fruta'
$ws.Range("G45").Value = 100109
$ws.Range('H45').Value = 'Uva'
$ws.Range("I45").Value = 100109001
$ws.Range('J45').Value = 'Uva'
$ws.Range('K45').Value = 'Superior Seedless'
$ws.Range('L45').Value = 'Primera'
$ws.Range("M45").Value = 400
$ws.Range("N45").Value = 11500
$ws.Range("O45").Value = 12000
$ws.Range("P45").Value = 11750
$ws.Range('Q45').Value = '$/caja 15 kilos'
$ws.Range('R45').Value = 'Provincia de Limarí'
$ws.Range("S45").Value = 783
$ws.Range("T45").Value = 15
